$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column K (11): widen the column to fit the new date format ---
$ws.Columns.Item(11).ColumnWidth = 16.6

# --- Re-format the Date column ---
# Header cell keeps a "dd-mm-yy" style
$ws.Range("K1").NumberFormat = "dd-mm-yy;@"

# The two data rows get reformatted to "yyyy-mm-dd" and bumped to the same date
$ws.Range("K2").NumberFormat = "yyyy-mm-dd;@"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd;@"
$ws.Range("K2").Value2 = 44587
$ws.Range("K3").Value2 = 44587

# The trailing blank cells in the column are no longer formatted/used
$ws.Range("K4").Clear()
$ws.Range("K5").Clear()

# --- Scroll the sheet so column K (and beyond) is visible, select K3 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$ws.Range("K3").Select()
